$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H6").Value = 1800.125
$ws.Range("I6").Value = 568
$ws.Range("J6").Value = 2539.4
$ws.Range("K6").Value = 1704
$ws.Range("L6").Value = 7618.200000000001
$ws.Range("M6").Value = -1592
$ws.Range("N6").Value = -7842.200000000001
$ws.Range("H34").Value = 5000
$ws.Range("I34").Value = 5000
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 5000
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -4797
$ws.Range("N34").ClearContents()
$ws.Range("H36").Value = 5000
$ws.Range("I36").Value = 5000
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 5000
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -4285
$ws.Range("N36").ClearContents()
$ws.Range("H55").Value = 296.66666
$ws.Range("I55").Value = 231.42857
$ws.Range("J55").Value = 525
$ws.Range("K55").Value = 231.42857
$ws.Range("L55").Value = 525
$ws.Range("M55").Value = -17.42857000000001
$ws.Range("N55").Value = -953
$ws.Range("H70").Value = 2639.6
$ws.Range("I70").Value = 2639.6
$ws.Range("K70").Value = 7918.799999999999
$ws.Range("M70").Value = -7648.799999999999
$ws.Range("H73").Value = 2639.6
$ws.Range("I73").Value = 2639.6
$ws.Range("K73").Value = 7918.799999999999
$ws.Range("M73").Value = -6982.799999999999
$ws.Range("H98").Value = 1794.6
$ws.Range("I98").Value = 2329.6667
$ws.Range("K98").Value = 2329.6667
$ws.Range("M98").Value = -831.6667000000002
$ws.Range("H101").Value = 2816.1667
$ws.Range("I101").Value = 4300
$ws.Range("J101").Value = 1332.3334
$ws.Range("K101").Value = 12900
$ws.Range("L101").Value = 3997.0002
$ws.Range("M101").Value = -11278
$ws.Range("N101").Value = -7241.0002
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H106").Value = 2728.5
$ws.Range("I106").Value = 2970
$ws.Range("J106").Value = 2004
$ws.Range("K106").Value = 2970
$ws.Range("L106").Value = 2004
$ws.Range("M106").Value = -2339
$ws.Range("N106").Value = -3266
$ws.Range("H122").Value = 1794.6
$ws.Range("I122").Value = 2329.6667
$ws.Range("K122").Value = 6989.000100000001
$ws.Range("M122").Value = -4539.000100000001
$ws.Range("H132").Value = 2261.16
$ws.Range("I132").Value = 1926.45
$ws.Range("K132").Value = 5779.35
$ws.Range("M132").Value = -3249.35
$ws.Range("H135").Value = 1287.091
$ws.Range("I135").Value = 462.1111
$ws.Range("J135").Value = 4999.5
$ws.Range("K135").Value = 4158.9999
$ws.Range("L135").Value = 44995.5
$ws.Range("M135").Value = -1623.9999
$ws.Range("N135").Value = -50065.5
$ws.Range("H137").Value = 1334.1666
$ws.Range("I137").Value = 1321.5385
$ws.Range("J137").Value = 1349.091
$ws.Range("K137").Value = 3964.6155
$ws.Range("L137").Value = 4047.273
$ws.Range("M137").Value = -1414.6155
$ws.Range("N137").Value = -9147.272999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 56592.332
$ws.Range("J24").Value = 56592.332
$ws.Range("L24").Value = 56592.332
$ws.Range("N24").Value = -57340.332
$ws.Range("H35").Value = 11750
$ws.Range("I35").Value = 3500
$ws.Range("J35").Value = 20000
$ws.Range("K35").Value = 3500
$ws.Range("L35").Value = 20000
$ws.Range("M35").Value = -3094
$ws.Range("N35").Value = -20812
$ws.Range("H74").Value = 802.6
$ws.Range("I74").Value = 802.6
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 802.6
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 71.39999999999998
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 802.6
$ws.Range("I77").Value = 802.6
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 4013
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 355
$ws.Range("N77").ClearContents()
$ws.Range("H95").Value = 72402.664
$ws.Range("J95").Value = 72402.664
$ws.Range("L95").Value = 72402.664
$ws.Range("N95").Value = -77894.664
$ws.Range("H96").Value = 28165.666
$ws.Range("J96").Value = 28165.666
$ws.Range("L96").Value = 28165.666
$ws.Range("N96").Value = -33657.666
$ws.Range("H97").Value = 2036.6666
$ws.Range("I97").Value = 3000
$ws.Range("J97").Value = 1844
$ws.Range("K97").Value = 3000
$ws.Range("L97").Value = 1844
$ws.Range("M97").Value = -2504
$ws.Range("N97").Value = -2836
$ws.Range("H100").Value = 56592.332
$ws.Range("J100").Value = 56592.332
$ws.Range("L100").Value = 56592.332
$ws.Range("N100").Value = -58756.332
$ws.Range("H101").Value = 50000
$ws.Range("J101").Value = 50000
$ws.Range("L101").Value = 50000
$ws.Range("N101").Value = -56490
$ws.Range("H132").Value = 977
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 9000
$ws.Range("J81").Value = 9000
$ws.Range("L81").Value = 9000
$ws.Range("N81").Value = -11122
$ws.Range("H84").Value = 9000
$ws.Range("J84").Value = 9000
$ws.Range("L84").Value = 27000
$ws.Range("N84").Value = -37608
$ws.Range("H94").Value = 3282.5715
$ws.Range("I94").Value = 999.5
$ws.Range("J94").Value = 4195.8
$ws.Range("K94").Value = 999.5
$ws.Range("L94").Value = 4195.8
$ws.Range("M94").Value = -548.5
$ws.Range("N94").Value = -5097.8
$ws.Range("H100").Value = 30643
$ws.Range("J100").Value = 30643
$ws.Range("L100").Value = 30643
$ws.Range("N100").Value = -32807

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 5000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 5000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 5000
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -5224

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4218.625
$ws.Range("I102").Value = 3541.5
$ws.Range("J102").Value = 6250
$ws.Range("K102").Value = 3541.5
$ws.Range("L102").Value = 6250
$ws.Range("M102").Value = -1919.5
$ws.Range("N102").Value = -9494

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 3515
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H93").Value = 867.1818
$ws.Range("I93").Value = 827.1429000000001
$ws.Range("K93").Value = 827.1429000000001
$ws.Range("M93").Value = 420.8570999999999
$ws.Range("H101").Value = 20362.666
$ws.Range("J101").Value = 20362.666
$ws.Range("L101").Value = 20362.666
$ws.Range("N101").Value = -26852.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 17275.25
$ws.Range("J103").Value = 17275.25
$ws.Range("L103").Value = 17275.25
$ws.Range("N103").Value = -19619.25
$ws.Range("H122").Value = 2240.8572
$ws.Range("I122").Value = 1737.4
$ws.Range("K122").Value = 5212.200000000001
$ws.Range("M122").Value = -2762.200000000001
$ws.Range("H136").Value = 674.96295
$ws.Range("I136").Value = 723.875
$ws.Range("J136").Value = 283.66666
$ws.Range("K136").Value = 2171.625
$ws.Range("L136").Value = 850.9999799999999
$ws.Range("M136").Value = 378.375
$ws.Range("N136").Value = -5950.99998
